# Scheduled-runner refresh of the Kujata leve-profit sheets.
# Columns H:N on each sheet hold point-in-time market data
# (currentAveragePrice*, LevePrice*, LeveProfit*) pulled live from the
# market board - this just writes the latest pull's numbers over the
# previous snapshot, row by row, on each affected leve sheet.
# Cells that the new pull didn't produce a figure for are cleared
# (their <c> element drops out of the row) rather than zeroed.

$wb = $excel.ActiveWorkbook

# --- ALC ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# r19 - Unbreak My Heart / Roof Tile
$ws.Range("H19").Value = 848
$ws.Range("I19").Value = 761.4286
$ws.Range("J19").Value = 894.61536
$ws.Range("K19").Value = 761.4286
$ws.Range("L19").Value = 894.61536
$ws.Range("M19").Value = -586.4286
$ws.Range("N19").Value = -1244.61536

# r113 - Amaro Kart / Starch Glue
$ws.Range("H113").Value = 3762.2
$ws.Range("I113").Value = 3502.5
$ws.Range("J113").Value = 3935.3333
$ws.Range("K113").Value = 3502.5
$ws.Range("L113").Value = 3935.3333
$ws.Range("M113").Value = -248.5
$ws.Range("N113").Value = -10443.3333

# r116 - Growing Up / Growth Formula Kappa
$ws.Range("H116").Value = 3143.75
$ws.Range("I116").Value = 2120
$ws.Range("J116").Value = 4167.5
$ws.Range("K116").Value = 2120
$ws.Range("L116").Value = 4167.5
$ws.Range("M116").Value = 1322
$ws.Range("N116").Value = -11051.5

# r137 - Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 1551.0555
$ws.Range("I137").Value = 1043
$ws.Range("K137").Value = 3129
$ws.Range("M137").Value = -579

# r138 - All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 1506.79
$ws.Range("I138").Value = 902.2174
$ws.Range("J138").Value = 1687.3766
$ws.Range("K138").Value = 2706.6522
$ws.Range("L138").Value = 5062.129800000001
$ws.Range("M138").Value = 2433.3478
$ws.Range("N138").Value = -15342.1298

# --- ARM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# r2 - Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 801.6923
$ws.Range("I2").Value = 619.9
$ws.Range("J2").Value = 1407.6666
$ws.Range("K2").Value = 619.9
$ws.Range("L2").Value = 1407.6666
$ws.Range("M2").Value = -506.9
$ws.Range("N2").Value = -1633.6666

# r12 - Strait Ain't the Gate / Bronze Scutum (no M figure this pull)
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

# r74 - As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 936.4286
$ws.Range("I74").Value = 753.6111
$ws.Range("J74").Value = 2033.3334
$ws.Range("K74").Value = 753.6111
$ws.Range("L74").Value = 2033.3334
$ws.Range("M74").Value = 120.3889
$ws.Range("N74").Value = -3781.3334

# r77 - Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 936.4286
$ws.Range("I77").Value = 753.6111
$ws.Range("J77").Value = 2033.3334
$ws.Range("K77").Value = 3768.0555
$ws.Range("L77").Value = 10166.667
$ws.Range("M77").Value = 599.9445000000001
$ws.Range("N77").Value = -18902.667

# r101 - Art Imitates Life / Doman Steel Tabard of Fending (no N figure this pull)
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

# r102 - Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Range("H102").Value = 11906918
$ws.Range("I102").Value = 15153445
$ws.Range("J102").Value = 2987.3333
$ws.Range("K102").Value = 15153445
$ws.Range("L102").Value = 2987.3333
$ws.Range("M102").Value = -15151823
$ws.Range("N102").Value = -6231.3333

# r116 - No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 801.6923
$ws.Range("I116").Value = 619.9
$ws.Range("J116").Value = 1407.6666
$ws.Range("K116").Value = 619.9
$ws.Range("L116").Value = 1407.6666
$ws.Range("M116").Value = 1674.1
$ws.Range("N116").Value = -5995.6666

# --- BSM -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# r3 - Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 801.6923
$ws.Range("I3").Value = 619.9
$ws.Range("J3").Value = 1407.6666
$ws.Range("K3").Value = 619.9
$ws.Range("L3").Value = 1407.6666
$ws.Range("M3").Value = -505.9
$ws.Range("N3").Value = -1635.6666

# --- CRP -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# r16 - Raise the Roof / Ash Lumber
$ws.Range("H16").Value = 50001390
$ws.Range("J16").Value = 1400
$ws.Range("L16").Value = 1400
$ws.Range("N16").Value = -1974

# r107 - Built to Last / White Oak Lumber
$ws.Range("H107").Value = 561.8570999999999
$ws.Range("I107").Value = 455.5
$ws.Range("K107").Value = 455.5
$ws.Range("M107").Value = 1464.5

# r113 - Patient Patients / White Ash Lumber
$ws.Range("H113").Value = 50001390
$ws.Range("J113").Value = 1400
$ws.Range("L113").Value = 1400
$ws.Range("N113").Value = -5740

# --- CUL -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# r131 - The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 21742768
$ws.Range("I131").Value = 200000260
$ws.Range("J131").Value = 4048.634
$ws.Range("K131").Value = 600000780
$ws.Range("L131").Value = 12145.902
$ws.Range("M131").Value = -599995740
$ws.Range("N131").Value = -22225.902

# --- GSM -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

# r92 - Play It by Ear / Triphane Earrings of Healing (no N figure this pull)
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# r97 - If I'd a Koppranickel for Every Time... / Koppranickel Ingot (N now present)
$ws.Range("H97").Value = 816.125
$ws.Range("I97").Value = 804.1429000000001
$ws.Range("J97").Value = 900
$ws.Range("K97").Value = 804.1429000000001
$ws.Range("L97").Value = 900
$ws.Range("M97").Value = -308.1429000000001
$ws.Range("N97").Value = -1892

# r102 - Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 1172.4117
$ws.Range("I102").Value = 966.5
$ws.Range("J102").Value = 2133.3333
$ws.Range("K102").Value = 966.5
$ws.Range("L102").Value = 2133.3333
$ws.Range("M102").Value = 655.5
$ws.Range("N102").Value = -5377.3333

# --- LTW -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# r61 - Spelling Me Softly / Raptor Leather (no M figure this pull)
$ws.Range("H61").Value = 1503.3334
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 1503.3334
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 1503.3334
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -1907.3334

# r68 - You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 2018.2174
$ws.Range("I68").Value = 1961
$ws.Range("K68").Value = 1961
$ws.Range("M68").Value = -1212

# r71 - They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 2018.2174
$ws.Range("I71").Value = 1961
$ws.Range("K71").Value = 9805
$ws.Range("M71").Value = -6061

# r113 - Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 1503.3334
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1503.3334
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1503.3334
$ws.Range("M113").Value = 1503.3334
$ws.Range("N113").Value = -5843.3334

# --- WVR -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# r100 - Of Great Import / Kudzu Thread
$ws.Range("H100").Value = 114.454544
$ws.Range("I100").Value = 96.5
$ws.Range("J100").Value = 136
$ws.Range("K100").Value = 193
$ws.Range("L100").Value = 272
$ws.Range("M100").Value = 348
$ws.Range("N100").Value = -1354

# r107 - Flax Wax / Bright Linen Yarn
$ws.Range("H107").Value = 508.2857
$ws.Range("I107").Value = 414.5
$ws.Range("K107").Value = 1243.5
$ws.Range("M107").Value = 676.5

# r126 - A Polished Purchase / Snow Linen
$ws.Range("H126").Value = 58825210
$ws.Range("I126").Value = 76924360
$ws.Range("K126").Value = 230773080
$ws.Range("M126").Value = -230770610

# r132 - Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 3218.5
$ws.Range("I132").Value = 2853.8333
$ws.Range("K132").Value = 8561.499899999999
$ws.Range("M132").Value = -6031.499899999999
